$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.168.33'
$ws.Range('E2').Value = '  -5.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.460.74'
$ws.Range('E3').Value = '  -7.09%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.52'
$ws.Range('E5').Value = '  -8.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.96'
$ws.Range('E6').Value = '  -3.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.456.92'
$ws.Range('E7').Value = '  -7.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  -6.29%  '
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.646'
$ws.Range('E10').Value = '  -10.85%  '
$ws.Range('E11').Value = '  -12.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.43'
$ws.Range('E12').Value = '  -13.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('E13').Value = '  -14.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.46'
$ws.Range('E14').Value = '  -11.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.035.30'
$ws.Range('E15').Value = '  -6.53%  '
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.489.33'
$ws.Range('E17').Value = '  -6.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.79'
$ws.Range('E18').Value = '  -8.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '65.040.89'
$ws.Range('E19').Value = '  -5.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.74'
$ws.Range('E20').Value = '  -9.46%  '
$ws.Range('E21').Value = '  -10.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '380.76'
$ws.Range('E22').Value = '  -7.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.15'
$ws.Range('E23').Value = '  -9.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.05'
$ws.Range('E24').Value = '  -7.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.68'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.81'
$ws.Range('E26').Value = '  -8.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.00'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.94'
$ws.Range('E28').Value = '  -7.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.39'
$ws.Range('E29').Value = '  -11.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.54'
$ws.Range('E30').Value = '  -11.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.13'
$ws.Range('E31').Value = '  -9.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.96'
$ws.Range('E32').Value = '  -8.32%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '604.93'
$ws.Range('E33').Value = '  -4.40%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.80'
$ws.Range('E34').Value = '  -7.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '62.34'
$ws.Range('E35').Value = '  -5.22%  '
$ws.Range('E36').Value = '  -10.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '40.35'
$ws.Range('E37').Value = '  -12.73%  '
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.386'
$ws.Range('E39').Value = '  -6.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0727'
$ws.Range('E40').Value = '  -12.44%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.128'
$ws.Range('E42').Value = '  -9.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.932.48'
$ws.Range('E43').Value = '  +1.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.71'
$ws.Range('E44').Value = '  -11.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.40'
$ws.Range('E45').Value = '  -8.65%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.10'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0393'
$ws.Range('E47').Value = '  -11.99%  '
$ws.Range('E48').Value = '  -9.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '137.14'
$ws.Range('E49').Value = '  -4.91%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.47'
$ws.Range('E50').Value = '  -10.58%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.16'
$ws.Range('E51').Value = '  -11.16%  '
